# Fixed StudyComb for Faceted Filters ICDC
#
# The "startup" sheet drives a set of Neo4j "StatQuery" cypher queries used
# to compute the header counts (files/samples/cases/studies) for each tab
# (CasesTab / SamplesTab / FilesTab). The old StatQuery text used a stale
# multi-MATCH pattern; replace it with the corrected query on every row
# that references it (column C, rows 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newStatQuery = "MATCH (demo:demographic)`nWHERE demo.breed IN [`"Akita`"]`nMATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`nOPTIONAL MATCH (c)<-[*]-(samp:sample)`nOPTIONAL MATCH (c)<-[*]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files, `n`tcount(DISTINCT(samp)) as number_of_sample, `n`tcount(DISTINCT(c)) as number_of_cases, `n`tcount(DISTINCT(s)) as number_of_study"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Reflect the author's last-saved cursor position (row 4 scrolled into view,
# B4 selected).
$ws.Activate()
$ws.Range("B4").Select()
